$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 6890
$ws.Range("F8").Value = 3439
$ws.Range("F9").Value = 3439
$ws.Range("F19").Value = 80
$ws.Range("F26").Value = 387
$ws.Range("F33").Value = 2030
$ws.Range("F38").Value = 3913
$ws.Range("F39").Value = 366
$ws.Range("F43").Value = 650
$ws.Range("F44").Value = 72
$ws.Range("F45").Value = 1506
$ws.Range("F46").Value = 253
$ws.Range("F48").Value = 579
$ws.Range("F49").Value = 669

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 404

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 148

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 148
$ws.Range("F10").Value = 6890
$ws.Range("F11").Value = 3439
$ws.Range("F12").Value = 3439
$ws.Range("F20").Value = 80
$ws.Range("F30").Value = 387
$ws.Range("F35").Value = 2030
$ws.Range("F40").Value = 3913
$ws.Range("F41").Value = 366
$ws.Range("F46").Value = 72
$ws.Range("F47").Value = 1506
$ws.Range("F48").Value = 253
$ws.Range("F50").Value = 669
